$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opv_results")

# Rows 10-17: full R_mean/RMSE_mean/R_std/RMSE_std/num_of_data block
$data1 = @(
    @([double]"0.61635443823195002", [double]"7.1095427005393896E-2", [double]"0.16054202616214699", [double]"2.4363607168197601E-2", 556),
    @([double]"0.62607187578574197", [double]"8.0043212497972194E-2", [double]"0.15966604650020599", [double]"2.56536491215229E-2", 556),
    @([double]"0.63635914334342403", [double]"6.5593093645026704E-2", [double]"0.156994849443435", [double]"2.34825499355793E-2", 556),
    @([double]"0.49433722894073201", [double]"7.7384885219848307E-2", [double]"0.176036447286605", [double]"2.41741314530372E-2", 556),
    @([double]"0.68772552132417297", [double]"8.1548127360052694E-2", [double]"0.14705356955528201", [double]"2.1178727969527199E-2", 556),
    @([double]"0.67575033016347297", [double]"7.3512473110954907E-2", [double]"0.148529052734375", [double]"1.48715740069746E-2", 556),
    @([double]"0.652411369232893", [double]"7.5018817297932902E-2", [double]"0.154740735888481", [double]"2.17382330447435E-2", 556),
    @([double]"0.71993180364388798", [double]"6.6144075928657203E-2", [double]"0.139440968632698", [double]"1.83883681893348E-2", 556)
)

$r = 10
foreach ($row in $data1) {
    $ws.Range("E$r").Value2 = $row[0]
    $ws.Range("F$r").Value2 = $row[1]
    $ws.Range("G$r").Value2 = $row[2]
    $ws.Range("H$r").Value2 = $row[3]
    $ws.Range("I$r").Value2 = $row[4]
    $r++
}

# Rows 42-49: full R_mean/RMSE_mean/R_std/RMSE_std/num_of_data block
$data2 = @(
    @([double]"0.63911453853812805", [double]"3.2522503023203503E-2", [double]"0.15611515512705701", [double]"1.21339535146861E-2", 447),
    @([double]"0.65682806789106596", [double]"2.9563898247692299E-2", [double]"0.15352006748400501", [double]"1.16529925882198E-2", 447),
    @([double]"0.651121451979078", [double]"1.5348962892391E-2", [double]"0.153760794688625", [double]"9.6016212375992799E-3", 447),
    @([double]"0.56312415504781199", [double]"3.8914309861472297E-2", [double]"0.17100185075917501", [double]"1.07003483878228E-2", 447),
    @([double]"0.69665226596533603", [double]"2.24192081331595E-2", [double]"0.14508134475361101", [double]"1.0029892229590999E-2", 447),
    @([double]"0.695116116823148", [double]"1.3071220611402101E-2", [double]"0.145537560342584", [double]"9.7381610904078205E-3", 447),
    @([double]"0.58590629073124401", [double]"5.6501301209272298E-2", [double]"0.16881166342197401", [double]"1.5218371091253799E-2", 447),
    @([double]"0.74228318921537395", [double]"1.9902063646396101E-2", [double]"0.13650432155313699", [double]"1.0602340360635501E-2", 447)
)

$r = 42
foreach ($row in $data2) {
    $ws.Range("E$r").Value2 = $row[0]
    $ws.Range("F$r").Value2 = $row[1]
    $ws.Range("G$r").Value2 = $row[2]
    $ws.Range("H$r").Value2 = $row[3]
    $ws.Range("I$r").Value2 = $row[4]
    $r++
}

$ws.Range("K25").Select()
